# Added webarchives for Bealto OpenCL Sorting.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$xlPasteFormats = -4122

# Row 30's date cell carried a one-off style (numFmtId 14); normalize it to the
# same date style (numFmtId 164) already used by the rest of column A, by
# copying the formatting from A2 onto it.
$ws.Range("A2").Copy()
$ws.Range("A30").PasteSpecial($xlPasteFormats)
$ws.Range("A30").Value = 41179

# New entry: 2012-10-01 - "Moved Buffer deletions into download method to avoid OUT_OF_RESOURCES"
$ws.Range("A31").PasteSpecial($xlPasteFormats)
$ws.Range("A31").Value = 41183
$ws.Range("B31").Value = "Moved Buffer deletions into download method to avoid OUT_OF_RESOURCES"

# Trailing blank, date-formatted rows reserved for future entries.
for ($r = 32; $r -le 39; $r++) {
    $ws.Cells.Item($r, 1).PasteSpecial($xlPasteFormats)
}
$excel.CutCopyMode = $false

$excel.Application.GoTo($ws.Range("B32"))
